# HORIZONTINA.xlsx update:
#  - rename "Paineis DARQ"            -> "PAINEIS DARQ"
#  - rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#  - remove "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Paineis DARQ")
$ws1.Name = "PAINEIS DARQ"

$ws5 = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$ws5.Name = "RECOLHIMENTO X ELIMINAÇÃO"

$ws7 = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$ws7.Delete() | Out-Null
